# Generate Report for Handoff
# The 0ba3c325-... row has moved from "Handed back: in sync with en-US"
# to "Ready for handoff" after a new handoff was generated. Update the
# Status, Latest Handoff Datetime and Error Detail columns (and the
# Overview roll-up sheet) accordingly.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1308ade67f2860791d17a79989cf4f22481a7356/e2e/0ba3c325-b4ba-4a0f-bbfa-54202c9eb50a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f400db10382b591d6427d0f9738e6554168f717a/e2e/0ba3c325-b4ba-4a0f-bbfa-54202c9eb50a.md."

# Excel's COM bridge stores ColumnWidth with a fixed +5/6 character offset
# versus the raw OOXML <col width=".."/> attribute, so back it out here to
# land on an exact width of 40 in the saved file.
$targetColWidth = 40 - (5 / 6)

# --- zh-cn sheet: row 3 is the 0ba3c325 file ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("H3").Value = "2016-08-27 10:48:12"
$wsZh.Range("P3").Value = $errorDetail
$wsZh.Columns.Item(16).ColumnWidth = $targetColWidth

# --- de-de sheet: row 3 is the 0ba3c325 file ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("H3").Value = "2016-08-27 10:48:17"
$wsDe.Range("P3").Value = $errorDetail
$wsDe.Columns.Item(16).ColumnWidth = $targetColWidth

# --- Overview sheet: row 3 is the 0ba3c325 file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-27 10:48:17"
